$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "ICS2121-1",
    "IMM3800-1",
    "IIC2733-1",
    "IEE2213-1",
    "IIQ3343-1",
    "IMT3150-1",
    "ICE2604-1",
    "IEE2463-1",
    "ICH3364-1",
    "ICM2803-1",
    "IIQ2673-1",
    "ICC2105-1",
    "IIC2733-2",
    "ICC3543-1",
    "ICE3613-1",
    "IEE3234-1",
    "ICS2563-1",
    "ICH2214-1",
    "ICE3124-1",
    "IEE2123-1",
    "ICS3413-1",
    "IMT3800-1",
    "ICE2028-1",
    "ICS3582-1",
    "IIC3800-1",
    "ICM3762-1",
    "IMM2003-1",
    "ICS2123-2",
    "IIC2713-3",
    "ICS2023-1",
    "ICM3243-1",
    "ICT3464-1",
    "ICC2304-1",
    "ICH2304-2",
    "IEE2713-1",
    "IMM3323-1",
    "IIC2713-1",
    "ICH1005-1",
    "ICH3600-1",
    "ICE2623-1",
    "IIC2713-2",
    "ICS3811-1",
    "ICT3623-1",
    "ICS2523-4",
    "IEE3732-1",
    "ICH3222-1",
    "IIC2213-1",
    "ICE2633-1",
    "IMM2033-1",
    "IIC3103-1",
    "IMT2111-1",
    "ICM3251-1",
    "ICS3762-1",
    "ICE2006-1",
    "ICC3264-1",
    "ICE2114-1",
    "IIQ3643-1",
    "IEE3373-1",
    "IIC3743-1",
    "IIC2613-1",
    "IIC2764-1",
    "ICM2313-1",
    "IMM2213-1",
    "IIC3113-1",
    "ICT2303-1",
    "ICT3523-1",
    "ICC3253-1",
    "ICM2203-1",
    "IEE2343-1",
    "ICS2123-3",
    "IEE2513-1",
    "IIC3113-2",
    "ICE3653-1",
    "ICH3350-1",
    "IEE2103-1",
    "ICE2313-1",
    "ICH2124-1",
    "IEE2613-1",
    "IEE2413-1",
    "IIC2333-1",
    "ICC3214-1",
    "IIC3143-1",
    "IIQ2133-1",
    "ICC2514-1",
    "ICE3753-1",
    "IMM2053-1",
    "ICS2523-3",
    "ICM2022-1",
    "ICE2533-1",
    "IIC3757-1",
    "ICE2703-1",
    "IMM2043-1",
    "ICH2114-1",
    "IMM3313-1",
    "ICE3663-1",
    "ICM2403-1",
    "IIC3242-1",
    "ICC1001-1",
    "IIC2133-1",
    "IIQ2043-1",
    "ICE3513-1",
    "ICC3434-1",
    "IIC2133-2",
    "ICM2213-1",
    "IEE2113-1",
    "ICT3283-1",
    "ICH2204-1",
    "IIQ2303-1",
    "ICM2333-1",
    "ICS3723-1",
    "ICH2304-1",
    "ICC2204-1",
    "ICE2020-1",
    "ICS3313-1",
    "ICE3443-1",
    "ICS2123-1",
    "ICH3374-1",
    "ICM3235-1",
    "ICS3151-1",
    "IMM2013-1",
    "ICH2574-1",
    "ICH3532-1",
    "ICE3233-1",
    "ICM1001-1",
    "ICM2223-1",
    "IIQ2003-1",
    "ICC3124-1",
    "IIC1001-1",
    "IMM1003-1",
    "ICT2233-1",
    "IIC2143-2",
    "IMT3130-1",
    "IIC2143-1",
    "ICM2413-1",
    "ICE3413-1",
    "ICS2563-2",
    "ICT3435-1",
    "ICM2003-1",
    "ICM2028-1",
    "IIQ2663-1",
    "A3_MS1 - Macroseccion 1",
    "B0_MS1 - Macroseccion 1",
    "B2_MS1 - Macroseccion 1",
    "C1_MS1 - Macroseccion 1",
    "I5_MS1 - Macroseccion 1",
    "I6_MS1 - Macroseccion 1",
    "ICE1514_Coordinado - Macroseccion",
    "ICH1104_Coordinado - Macroseccion",
    "ICM2503_Coordinado - Macroseccion",
    "ICS1113_Coordinado - Macroseccion",
    "ICS1513_Coordinado - Macroseccion",
    "ICS2523_1_2_Sec_Coordinadas - Macroseccion",
    "ICS2613_Coordinado - Macroseccion",
    "ICS3213_Coordinado - Macroseccion",
    "ICS3313_2_3_Sec_Coordinadas - Macroseccion",
    "ICS3413_2_3_Sec_Coordinadas - Macroseccion",
    "IIC1103_Coordinado - Macroseccion",
    "IIC1253_Coordinado - Macroseccion",
    "IIC2343_Coordinado - Macroseccion",
    "IIC3745_Coordinado - Macroseccion",
    "IIQ1003_Coordinado - Macroseccion",
    "J7_MS1 - Macroseccion 1",
    "R2_MS1 - Macroseccion 1",
    "R3_MS1 - Macroseccion 1",
    "W0_MS1 - Macroseccion 1",
    "Z8_MS1 - Macroseccion 1",
    "EYP1113_Coordinado - Macroseccion",
    "FIS1514_Coordinado - Macroseccion",
    "FIS1523_Coordinado - Macroseccion",
    "FIS1533_Coordinado - Macroseccion",
    "MAT1203_Coordinado - Macroseccion",
    "MAT1207_Coordinado - Macroseccion",
    "MAT1610_Coordinado - Macroseccion",
    "MAT1620_Coordinado - Macroseccion",
    "MAT1630_Coordinado - Macroseccion",
    "MAT1640_Coordinado - Macroseccion",
    "QIM100E_Coordinado - Macroseccion"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

Write-Host "Done setting values"
Write-Host $ws.Range("A2").Value2
Write-Host $ws.Range("A178").Value2